# "updating the team diary"
# - Sheet1 keeps all its data but gets a few cosmetic view/row-height tweaks
#   (selection range, a couple of row heights) from being re-opened/edited.
# - A new Sheet2 is added after Sheet1 and becomes the active sheet. It starts
#   as a duplicate of Sheet1 (same headers/styles/layout) then is repurposed
#   as the diary entry for the first team meeting: new members-present text,
#   a new meeting time/discussion in row 7, and the future placeholder
#   meeting rows (10,12,14,16,18) are wiped back to blank template rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create Sheet2 as a copy of Sheet1 (before Sheet1's own row heights
#     below change), placed right after it ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# --- Sheet1: minor view / row-height adjustments (content unchanged) ---
$ws1.Rows.Item(7).RowHeight = 85
$ws1.Rows.Item(10).RowHeight = 34
$ws1.Rows.Item(12).RowHeight = 34
$ws1.Range("A1:E23").Select()

# New group-members text for this diary entry (single combined cell, second
# member cell no longer used)
$ws2.Range("B2").Value = "Evan,Amrita,Phuong Mai"
$ws2.Range("C2").ClearContents()

# First team meeting: new date/time details + discussion notes
$ws2.Range("B7").Value = 0.42708333333333331
$ws2.Range("C7").Value = 0.45833333333333331
$ws2.Range("E7").Value = "First team meeting completed - disucussing the action plan to proceed further "
$ws2.Rows.Item(7).RowHeight = 34

# Clear out the placeholder future-meeting rows back to blank template rows
$ws2.Range("A10:E10").ClearContents()
$ws2.Range("A12:E12").ClearContents()
$ws2.Range("A14:E14").ClearContents()
$ws2.Range("A16:E16").ClearContents()
$ws2.Range("A18:E18").ClearContents()

# Leave selection parked where the author left it
$ws2.Range("B21").Select()
